# Updates the cryptos price/volume table (and fixes two rows whose
# ranking order flipped) to match the latest scrape.
#
# Note: several "Price" (column D) values look like numbers to Excel
# (e.g. "148.01"), but the sheet stores them as literal text so that
# trailing zeros / exact digit strings survive untouched. Prefixing
# such values with a leading apostrophe forces Excel to keep them as
# text instead of silently coercing them into doubles (which would
# turn "1.00" into "1", "80.60" into "80.6", etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.704.58"
$ws.Range("E2").Value = "  -2.50%  "

# Row 3
$ws.Range("D3").Value = "2.906.37"
$ws.Range("E3").Value = "  -3.67%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'586.75"
$ws.Range("E5").Value = "  -1.30%  "

# Row 6
$ws.Range("D6").Value = "'148.01"
$ws.Range("E6").Value = "  -0.30%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("E8").Value = "  -2.71%  "

# Row 9
$ws.Range("D9").Value = "2.905.41"
$ws.Range("E9").Value = "  -3.65%  "

# Row 10
$ws.Range("D10").Value = "'6.72"
$ws.Range("E10").Value = "  +5.19%  "

# Row 11
$ws.Range("E11").Value = "  -3.78%  "

# Row 12
$ws.Range("E12").Value = "  -2.29%  "

# Row 13
$ws.Range("D13").Value = "'0.0000224"
$ws.Range("E13").Value = "  -3.61%  "

# Row 14
$ws.Range("D14").Value = "'34.11"
$ws.Range("E14").Value = "  -0.87%  "

# Row 15
$ws.Range("E15").Value = "  +0.50%  "

# Row 16
$ws.Range("D16").Value = "3.388.95"
$ws.Range("E16").Value = "  -3.55%  "

# Row 17
$ws.Range("D17").Value = "'6.82"
$ws.Range("E17").Value = "  -2.59%  "

# Row 18
$ws.Range("D18").Value = "60.686.28"
$ws.Range("E18").Value = "  -2.45%  "

# Row 19
$ws.Range("D19").Value = "2.908.40"
$ws.Range("E19").Value = "  -3.62%  "

# Row 20
$ws.Range("D20").Value = "'427.36"
$ws.Range("E20").Value = "  -4.75%  "

# Row 21
$ws.Range("D21").Value = "'13.64"
$ws.Range("E21").Value = "  -4.01%  "

# Row 22
$ws.Range("D22").Value = "'0.671"
$ws.Range("E22").Value = "  -2.79%  "

# Row 23
$ws.Range("E23").Value = "  -3.93%  "

# Row 24 -> Litecoin
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'80.60"
$ws.Range("E24").Value = "  -2.17%  "

# Row 25 -> RenderToken
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "'11.08"
$ws.Range("E25").Value = "  +1.51%  "

# Row 26
$ws.Range("D26").Value = "'2.22"
$ws.Range("E26").Value = "  -1.20%  "

# Row 27
$ws.Range("D27").Value = "'11.84"
$ws.Range("E27").Value = "  -1.66%  "

# Row 28
$ws.Range("E28").Value = "  +0.05%  "

# Row 29 -> FirstDigitalUSD
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.11%  "

# Row 30 -> NEARProtocol
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.24"
$ws.Range("E30").Value = "  +0.72%  "

# Row 31 -> PancakeSwap
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.62"
$ws.Range("E31").Value = "  -2.99%  "

# Row 32 -> ImmutableX
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.18"
$ws.Range("E32").Value = "  +2.32%  "

# Row 33
$ws.Range("D33").Value = "'26.52"
$ws.Range("E33").Value = "  -3.69%  "

# Row 34
$ws.Range("E34").Value = "  -3.21%  "

# Row 35
$ws.Range("E35").Value = "  -1.43%  "

# Row 36
$ws.Range("E36").Value = "  -2.15%  "

# Row 37
$ws.Range("D37").Value = "'5.68"
$ws.Range("E37").Value = "  -2.83%  "

# Row 38
$ws.Range("D38").Value = "'2.03"
$ws.Range("E38").Value = "  -1.48%  "

# Row 39 -> dogwifhat
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'2.96"
$ws.Range("E39").Value = "  -0.64%  "

# Row 40 -> OKB
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'49.32"
$ws.Range("E40").Value = "  -1.72%  "

# Row 41
$ws.Range("D41").Value = "'8.74"
$ws.Range("E41").Value = "  -3.86%  "

# Row 42
$ws.Range("E42").Value = "  -1.50%  "

# Row 43
$ws.Range("D43").Value = "'0.292"
$ws.Range("E43").Value = "  +2.05%  "

# Row 44
$ws.Range("D44").Value = "'41.59"

# Row 45
$ws.Range("D45").Value = "'0.0348"
$ws.Range("E45").Value = "  -1.09%  "

# Row 46
$ws.Range("D46").Value = "'370.71"
$ws.Range("E46").Value = "  -6.09%  "

# Row 47 -> Monero
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'133.39"
$ws.Range("E47").Value = "  -1.00%  "

# Row 48 -> Maker
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.655.14"
$ws.Range("E48").Value = "  -2.80%  "

# Row 50
$ws.Range("D50").Value = "'25.16"
$ws.Range("E50").Value = "  +5.94%  "

# Row 51
$ws.Range("E51").Value = "  -1.35%  "
